$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 17
$ws_ALC.Range("H17").Value = 388946.03
$ws_ALC.Range("I17").Value = 6952.5
$ws_ALC.Range("J17").Value = 503544.1
$ws_ALC.Range("K17").Value = 20857.5
$ws_ALC.Range("L17").Value = 1510632.3
$ws_ALC.Range("M17").Value = -20689.5
$ws_ALC.Range("N17").Value = -1510968.3

# ALC row 64
$ws_ALC.Range("H64").Value = 4243.5
$ws_ALC.Range("J64").Value = 4500
$ws_ALC.Range("L64").Value = 4500
$ws_ALC.Range("N64").Value = -4996

# ALC row 67
$ws_ALC.Range("H67").Value = 4243.5
$ws_ALC.Range("J67").Value = 4500
$ws_ALC.Range("L67").Value = 4500
$ws_ALC.Range("N67").Value = -6216

# ALC row 86
$ws_ALC.Range("H86").Value = 4704.095
$ws_ALC.Range("I86").Value = 3108.2424
$ws_ALC.Range("J86").Value = 10555.556
$ws_ALC.Range("K86").Value = 3108.2424
$ws_ALC.Range("L86").Value = 10555.556
$ws_ALC.Range("M86").Value = -1985.2424
$ws_ALC.Range("N86").Value = -12801.556

# ALC row 89
$ws_ALC.Range("H89").Value = 4704.095
$ws_ALC.Range("I89").Value = 3108.2424
$ws_ALC.Range("J89").Value = 10555.556
$ws_ALC.Range("K89").Value = 15541.212
$ws_ALC.Range("L89").Value = 52777.78
$ws_ALC.Range("M89").Value = -9925.212
$ws_ALC.Range("N89").Value = -64009.78

# ALC row 125
$ws_ALC.Range("H125").Value = 850
$ws_ALC.Range("I125").Value = 700
$ws_ALC.Range("J125").Value = 1000
$ws_ALC.Range("K125").Value = 6300
$ws_ALC.Range("L125").Value = 9000
$ws_ALC.Range("M125").Value = -3840
$ws_ALC.Range("N125").Value = -13920

# ARM row 6
$ws_ARM.Range("H6").Value = 0
$ws_ARM.Range("I6").Value = 0
$ws_ARM.Range("K6").Value = 0
$ws_ARM.Range("M6").ClearContents()

# ARM row 74
$ws_ARM.Range("H74").Value = 3395.9
$ws_ARM.Range("I74").Value = 2995.625
$ws_ARM.Range("J74").Value = 4997
$ws_ARM.Range("K74").Value = 2995.625
$ws_ARM.Range("L74").Value = 4997
$ws_ARM.Range("M74").Value = -2121.625
$ws_ARM.Range("N74").Value = -6745

# ARM row 77
$ws_ARM.Range("H77").Value = 3395.9
$ws_ARM.Range("I77").Value = 2995.625
$ws_ARM.Range("J77").Value = 4997
$ws_ARM.Range("K77").Value = 14978.125
$ws_ARM.Range("L77").Value = 24985
$ws_ARM.Range("M77").Value = -10610.125
$ws_ARM.Range("N77").Value = -33721

# ARM row 131
$ws_ARM.Range("H131").Value = 0
$ws_ARM.Range("J131").Value = 0
$ws_ARM.Range("L131").Value = 0
$ws_ARM.Range("N131").ClearContents()

# BSM row 134
$ws_BSM.Range("H134").Value = 1563.4706
$ws_BSM.Range("I134").Value = 1371.9333
$ws_BSM.Range("K134").Value = 4115.7999
$ws_BSM.Range("M134").Value = -1580.7999

# CRP row 6
$ws_CRP.Range("H6").Value = 3752097.8
$ws_CRP.Range("I6").Value = 5627147
$ws_CRP.Range("K6").Value = 5627147
$ws_CRP.Range("M6").Value = -5627034

# CRP row 7
$ws_CRP.Range("H7").Value = 79.3
$ws_CRP.Range("I7").Value = 77
$ws_CRP.Range("J7").Value = 100
$ws_CRP.Range("K7").Value = 77
$ws_CRP.Range("L7").Value = 100
$ws_CRP.Range("M7").Value = 36
$ws_CRP.Range("N7").Value = -326

# CRP row 22
$ws_CRP.Range("H22").Value = 2061.75
$ws_CRP.Range("I22").Value = 499.5
$ws_CRP.Range("K22").Value = 499.5
$ws_CRP.Range("M22").Value = -149.5

# CRP row 99
$ws_CRP.Range("H99").Value = 2624
$ws_CRP.Range("I99").Value = 2405
$ws_CRP.Range("K99").Value = 2405
$ws_CRP.Range("M99").Value = -907

# CRP row 126
$ws_CRP.Range("H126").Value = 2624
$ws_CRP.Range("I126").Value = 2405
$ws_CRP.Range("K126").Value = 7215
$ws_CRP.Range("M126").Value = -4745

# CUL row 14
$ws_CUL.Range("H14").Value = 4333.3335
$ws_CUL.Range("I14").Value = 4333.3335
$ws_CUL.Range("K14").Value = 13000.0005
$ws_CUL.Range("M14").Value = -12827.0005

# CUL row 97
$ws_CUL.Range("H97").Value = 1584.5555
$ws_CUL.Range("I97").Value = 1559.25
$ws_CUL.Range("K97").Value = 4677.75
$ws_CUL.Range("M97").Value = -4181.75

# CUL row 121
$ws_CUL.Range("H121").Value = 70311.72
$ws_CUL.Range("I121").Value = 1237.1666
$ws_CUL.Range("J121").Value = 104849
$ws_CUL.Range("K121").Value = 3711.4998
$ws_CUL.Range("L121").Value = 314547
$ws_CUL.Range("M121").Value = -2401.4998
$ws_CUL.Range("N121").Value = -317167

# CUL row 134
$ws_CUL.Range("H134").Value = 1870.8182
$ws_CUL.Range("I134").Value = 1870.8182
$ws_CUL.Range("K134").Value = 5612.4546
$ws_CUL.Range("M134").Value = -542.4546

# CUL row 139
$ws_CUL.Range("H139").Value = 66680584
$ws_CUL.Range("I139").Value = 76205670
$ws_CUL.Range("K139").Value = 228617010
$ws_CUL.Range("M139").Value = -228611870

# GSM row 70
$ws_GSM.Range("H70").Value = 8714.968999999999
$ws_GSM.Range("I70").Value = 8444.552
$ws_GSM.Range("K70").Value = 8444.552
$ws_GSM.Range("M70").Value = -8174.552

# GSM row 73
$ws_GSM.Range("H73").Value = 8714.968999999999
$ws_GSM.Range("I73").Value = 8444.552
$ws_GSM.Range("K73").Value = 8444.552
$ws_GSM.Range("M73").Value = -7508.552

# GSM row 80
$ws_GSM.Range("H80").Value = 2039.0416
$ws_GSM.Range("I80").Value = 2057.2144
$ws_GSM.Range("K80").Value = 2057.2144
$ws_GSM.Range("M80").Value = -1059.2144

# GSM row 83
$ws_GSM.Range("H83").Value = 2039.0416
$ws_GSM.Range("I83").Value = 2057.2144
$ws_GSM.Range("K83").Value = 10286.072
$ws_GSM.Range("M83").Value = -5294.072

# GSM row 102
$ws_GSM.Range("H102").Value = 1442.2162
$ws_GSM.Range("I102").Value = 701.96295
$ws_GSM.Range("K102").Value = 701.96295
$ws_GSM.Range("M102").Value = 920.03705

# GSM row 113
$ws_GSM.Range("H113").Value = 4121.3
$ws_GSM.Range("I113").Value = 1400
$ws_GSM.Range("K113").Value = 1400
$ws_GSM.Range("M113").Value = 770

# GSM row 132
$ws_GSM.Range("H132").Value = 3209.1025
$ws_GSM.Range("J132").Value = 4482
$ws_GSM.Range("L132").Value = 13446
$ws_GSM.Range("N132").Value = -18506

# LTW row 7
$ws_LTW.Range("H7").Value = 27781414
$ws_LTW.Range("I7").Value = 45457376
$ws_LTW.Range("J7").Value = 4900.5713
$ws_LTW.Range("K7").Value = 45457376
$ws_LTW.Range("L7").Value = 4900.5713
$ws_LTW.Range("M7").Value = -45457264
$ws_LTW.Range("N7").Value = -5124.5713

# LTW row 46
$ws_LTW.Range("H46").Value = 1627.4667
$ws_LTW.Range("I46").Value = 750
$ws_LTW.Range("J46").Value = 1690.1428
$ws_LTW.Range("K46").Value = 750
$ws_LTW.Range("L46").Value = 1690.1428
$ws_LTW.Range("M46").Value = -562
$ws_LTW.Range("N46").Value = -2066.1428

# LTW row 126
$ws_LTW.Range("H126").Value = 27781414
$ws_LTW.Range("I126").Value = 45457376
$ws_LTW.Range("J126").Value = 4900.5713
$ws_LTW.Range("K126").Value = 136372128
$ws_LTW.Range("L126").Value = 14701.7139
$ws_LTW.Range("M126").Value = -136369658
$ws_LTW.Range("N126").Value = -19641.7139

# WVR row 132
$ws_WVR.Range("H132").Value = 4993.84
$ws_WVR.Range("I132").Value = 5013.722
$ws_WVR.Range("K132").Value = 15041.166
$ws_WVR.Range("M132").Value = -12511.166
